# This script applies a batch market-data refresh to the Pandaemonium
# Profits workbook, updating the pulled-price / profit columns
# (H:N) on each job sheet to the values captured by the latest
# scheduled run, per Sheets/Pandaemonium_Profits.xlsx.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 455.66666
$ws.Range("I18").Value = 387.625
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 387.625
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = -103.625
$ws.Range("N18").Value = -1568
$ws.Range("H40").Value = 2428.4375
$ws.Range("I40").Value = 2864.3
$ws.Range("J40").Value = 1702
$ws.Range("K40").Value = 2864.3
$ws.Range("L40").Value = 1702
$ws.Range("M40").Value = -2689.3
$ws.Range("N40").Value = -2052
$ws.Range("H64").Value = 3939.12
$ws.Range("I64").Value = 3561.125
$ws.Range("K64").Value = 3561.125
$ws.Range("M64").Value = -3313.125
$ws.Range("H67").Value = 3939.12
$ws.Range("I67").Value = 3561.125
$ws.Range("K67").Value = 3561.125
$ws.Range("M67").Value = -2703.125
$ws.Range("H132").Value = 1595.2632
$ws.Range("I132").Value = 1246.5193
$ws.Range("J132").Value = 5222.2
$ws.Range("K132").Value = 3739.5579
$ws.Range("L132").Value = 15666.6
$ws.Range("M132").Value = -1209.5579
$ws.Range("N132").Value = -20726.6
$ws.Range("H138").Value = 1140918
$ws.Range("I138").Value = 2672.5264
$ws.Range("J138").Value = 1495453.4
$ws.Range("K138").Value = 8017.5792
$ws.Range("L138").Value = 4486360.199999999
$ws.Range("M138").Value = -2877.5792
$ws.Range("N138").Value = -4496640.199999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13137.581
$ws.Range("I32").Value = 10406.415
$ws.Range("J32").Value = 29221.111
$ws.Range("K32").Value = 10406.415
$ws.Range("L32").Value = 29221.111
$ws.Range("M32").Value = -10119.415
$ws.Range("N32").Value = -29795.111
$ws.Range("H37").Value = 17999
$ws.Range("J37").Value = 24666.666
$ws.Range("L37").Value = 24666.666
$ws.Range("N37").Value = -25212.666
$ws.Range("H45").Value = 1435.1818
$ws.Range("I45").Value = 1407.3448
$ws.Range("J45").Value = 1637
$ws.Range("K45").Value = 1407.3448
$ws.Range("L45").Value = 1637
$ws.Range("M45").Value = -1030.3448
$ws.Range("N45").Value = -2391
$ws.Range("H61").Value = 11023.963
$ws.Range("I61").Value = 6704.533
$ws.Range("K61").Value = 6704.533
$ws.Range("M61").Value = -6492.533
$ws.Range("H63").Value = 10616.308
$ws.Range("I63").Value = 2333.3333
$ws.Range("K63").Value = 2333.3333
$ws.Range("M63").Value = -1647.3333
$ws.Range("H66").Value = 10616.308
$ws.Range("I66").Value = 2333.3333
$ws.Range("K66").Value = 11666.6665
$ws.Range("M66").Value = -8234.666499999999
$ws.Range("H74").Value = 4138.886
$ws.Range("I74").Value = 2107.6206
$ws.Range("J74").Value = 8066
$ws.Range("K74").Value = 2107.6206
$ws.Range("L74").Value = 8066
$ws.Range("M74").Value = -1233.6206
$ws.Range("N74").Value = -9814
$ws.Range("H77").Value = 4138.886
$ws.Range("I77").Value = 2107.6206
$ws.Range("J77").Value = 8066
$ws.Range("K77").Value = 10538.103
$ws.Range("L77").Value = 40330
$ws.Range("M77").Value = -6170.103000000001
$ws.Range("N77").Value = -49066
$ws.Range("H122").Value = 15627665
$ws.Range("I122").Value = 2663.8
$ws.Range("K122").Value = 7991.400000000001
$ws.Range("M122").Value = -5541.400000000001
$ws.Range("H132").Value = 5215.659
$ws.Range("I132").Value = 2091.9443
$ws.Range("J132").Value = 7378.231
$ws.Range("K132").Value = 6275.8329
$ws.Range("L132").Value = 22134.693
$ws.Range("M132").Value = -3745.8329
$ws.Range("N132").Value = -27194.693
$ws.Range("H136").Value = 11023.963
$ws.Range("I136").Value = 6704.533
$ws.Range("K136").Value = 20113.599
$ws.Range("M136").Value = -17563.599

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 18100.8
$ws.Range("I7").Value = 750
$ws.Range("J7").Value = 29668
$ws.Range("K7").Value = 750
$ws.Range("L7").Value = 29668
$ws.Range("M7").Value = -637
$ws.Range("N7").Value = -29894
$ws.Range("H25").Value = 3034.1667
$ws.Range("I25").Value = 1841
$ws.Range("J25").Value = 9000
$ws.Range("K25").Value = 1841
$ws.Range("L25").Value = 9000
$ws.Range("M25").Value = -1606
$ws.Range("N25").Value = -9470
$ws.Range("H35").Value = 36074
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H132").Value = 62981.668
$ws.Range("J132").Value = 62981.668
$ws.Range("L132").Value = 62981.668
$ws.Range("N132").Value = -73101.66800000001
$ws.Range("H134").Value = 29100.486
$ws.Range("I134").Value = 3288.037
$ws.Range("K134").Value = 9864.110999999999
$ws.Range("M134").Value = -7329.110999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 84.5
$ws.Range("I7").Value = 79
$ws.Range("K7").Value = 79
$ws.Range("M7").Value = 34
$ws.Range("H31").Value = 3228.697
$ws.Range("I31").Value = 1119.6842
$ws.Range("J31").Value = 6090.9287
$ws.Range("K31").Value = 1119.6842
$ws.Range("L31").Value = 6090.9287
$ws.Range("M31").Value = -824.6841999999999
$ws.Range("N31").Value = -6680.9287
$ws.Range("H34").Value = 3228.697
$ws.Range("I34").Value = 1119.6842
$ws.Range("J34").Value = 6090.9287
$ws.Range("K34").Value = 1119.6842
$ws.Range("L34").Value = 6090.9287
$ws.Range("M34").Value = -917.6841999999999
$ws.Range("N34").Value = -6494.9287
$ws.Range("H58").Value = 2758390
$ws.Range("I58").Value = 5053421.5
$ws.Range("J58").Value = 4352
$ws.Range("K58").Value = 5053421.5
$ws.Range("L58").Value = 4352
$ws.Range("M58").Value = -5053218.5
$ws.Range("N58").Value = -4758
$ws.Range("H59").Value = 19700
$ws.Range("J59").Value = 19888.889
$ws.Range("L59").Value = 19888.889
$ws.Range("N59").Value = -22178.889
$ws.Range("H60").Value = 18925
$ws.Range("J60").Value = 27633.334
$ws.Range("L60").Value = 27633.334
$ws.Range("N60").Value = -28655.334
$ws.Range("H62").Value = 3067.3333
$ws.Range("I62").Value = 3080.8
$ws.Range("K62").Value = 3080.8
$ws.Range("M62").Value = -2456.8
$ws.Range("H65").Value = 3067.3333
$ws.Range("I65").Value = 3080.8
$ws.Range("K65").Value = 15404
$ws.Range("M65").Value = -12284
$ws.Range("H74").Value = 33314
$ws.Range("J74").Value = 33314
$ws.Range("L74").Value = 33314
$ws.Range("N74").Value = -35062
$ws.Range("H77").Value = 33314
$ws.Range("J77").Value = 33314
$ws.Range("L77").Value = 99942
$ws.Range("N77").Value = -108678
$ws.Range("H136").Value = 2758390
$ws.Range("I136").Value = 5053421.5
$ws.Range("J136").Value = 4352
$ws.Range("K136").Value = 15160264.5
$ws.Range("L136").Value = 13056
$ws.Range("M136").Value = -15157714.5
$ws.Range("N136").Value = -18156

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2566857.5
$ws.Range("I5").Value = 615
$ws.Range("J5").Value = 5752538
$ws.Range("K5").Value = 1845
$ws.Range("L5").Value = 17257614
$ws.Range("M5").Value = -1733
$ws.Range("N5").Value = -17257838
$ws.Range("H107").Value = 2365082
$ws.Range("J107").Value = 1207.6945
$ws.Range("L107").Value = 3623.0835
$ws.Range("N107").Value = -7463.083500000001
$ws.Range("H131").Value = 540.5599999999999
$ws.Range("I131").Value = 292.2549
$ws.Range("J131").Value = 799
$ws.Range("K131").Value = 876.7647000000001
$ws.Range("L131").Value = 2397
$ws.Range("M131").Value = 4163.2353
$ws.Range("N131").Value = -12477
$ws.Range("H132").Value = 1486.3914
$ws.Range("J132").Value = 2301
$ws.Range("L132").Value = 20709
$ws.Range("N132").Value = -25769
$ws.Range("H135").Value = 2566857.5
$ws.Range("I135").Value = 615
$ws.Range("J135").Value = 5752538
$ws.Range("K135").Value = 5535
$ws.Range("L135").Value = 51772842
$ws.Range("M135").Value = -3000
$ws.Range("N135").Value = -51777912
$ws.Range("H137").Value = 38466772
$ws.Range("I137").Value = 45458268
$ws.Range("J137").Value = 13532.5
$ws.Range("K137").Value = 136374804
$ws.Range("L137").Value = 40597.5
$ws.Range("M137").Value = -136369704
$ws.Range("N137").Value = -50797.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5954
$ws.Range("J80").Value = 4223.5293
$ws.Range("L80").Value = 4223.5293
$ws.Range("N80").Value = -6219.5293
$ws.Range("H83").Value = 5954
$ws.Range("J83").Value = 4223.5293
$ws.Range("L83").Value = 21117.6465
$ws.Range("N83").Value = -31101.6465
$ws.Range("H132").Value = 38048.094
$ws.Range("I132").Value = 70528.53
$ws.Range("K132").Value = 211585.59
$ws.Range("M132").Value = -209055.59

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1601
$ws.Range("I46").Value = 1601
$ws.Range("K46").Value = 1601
$ws.Range("M46").Value = -1413
$ws.Range("H68").Value = 1850.25
$ws.Range("J68").Value = 1960
$ws.Range("L68").Value = 1960
$ws.Range("N68").Value = -3458
$ws.Range("H71").Value = 1850.25
$ws.Range("J71").Value = 1960
$ws.Range("L71").Value = 9800
$ws.Range("N71").Value = -17288
$ws.Range("H136").Value = 7544.577
$ws.Range("I136").Value = 6535.8184
$ws.Range("J136").Value = 8284.333000000001
$ws.Range("K136").Value = 19607.4552
$ws.Range("L136").Value = 24852.999
$ws.Range("M136").Value = -17057.4552
$ws.Range("N136").Value = -29952.999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 29072
$ws.Range("J97").Value = 29072
$ws.Range("L97").Value = 29072
$ws.Range("N97").Value = -31054
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H132").Value = 1301.3541
$ws.Range("I132").Value = 763.0732
$ws.Range("J132").Value = 4454.143
$ws.Range("K132").Value = 2289.2196
$ws.Range("L132").Value = 13362.429
$ws.Range("M132").Value = 240.7803999999996
$ws.Range("N132").Value = -18422.429

Write-Output "Applied all cell updates."
